# This script applies the edit described by the commit:
# "added script for failing a tese case"
#
# Semantic changes:
# 1. Workbook view: the active tab moves from sheet2 (createLeadTest) to sheet1 (loginTest).
# 2. The shared string "Lightning Experience" is replaced by a new, separate string
#    "Lightning" used by loginTest!C2 (homepageTitle value).
# 3. loginTest's two hyperlinks (A2 / B2) are swapped, so they now point at the wrong
#    mailto targets (intentionally introduces a mismatch / failing test case).

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("loginTest")

# 1) Change the homepageTitle value used in loginTest from "Lightning Experience" to "Lightning"
$wsLogin.Range("C2").Value = "Lightning"

# 2) Swap the hyperlinks on loginTest between A2 and B2.
$wsLogin.Hyperlinks.Delete()
$wsLogin.Hyperlinks.Add($wsLogin.Range("A2"), "mailto:Umang@8880")
$wsLogin.Hyperlinks.Add($wsLogin.Range("B2"), "mailto:umang8880@up.com")

# Re-adding hyperlinks resets the cell style to a fresh "Hyperlink"-like style
# instead of reusing the workbook's existing named style; restore the
# original shared style explicitly so cell formatting is unchanged.
$wsLogin.Range("A2").Style = "Hyperlink"
$wsLogin.Range("B2").Style = "Hyperlink"

# 3) Switch which sheet is the active/selected tab: loginTest becomes selected,
#    createLeadTest is no longer the selected tab.
$wsLogin.Activate()
$wsLogin.Select()
